$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update response text (column B) for rows 2-9.
# B5 and B7 become purely-numeric-looking strings ("450" / "100"); force
# text formatting first so Excel keeps them as text cells instead of
# auto-converting them to numbers.
$ws.Range("B2").Value = "You can have up to 250 curve shades per plot."
$ws.Range("B3").Value = "Unlimited"
$ws.Range("B4").Value = "The name ""Hydrocarbon bearing zone highlighted"" exceeds the maximum allowed length of 20 characters for curve shade names."

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "450"

$ws.Range("B6").Value = "You have reached the maximum limit of 20000 modifier types per plot."

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "100"

$ws.Range("B8").Value = "You are allowed to define 23 tadpole definitions."
$ws.Range("B9").Value = "You have reached the maximum of 50 layouts per ODF file."

# Update running time (column C) values for rows 2-9
$ws.Range("C2").Value = 19.28
$ws.Range("C3").Value = 16.37
$ws.Range("C4").Value = 43.64
$ws.Range("C5").Value = 15.93
$ws.Range("C6").Value = 14.52
$ws.Range("C7").Value = 20.62
$ws.Range("C8").Value = 37.18
$ws.Range("C9").Value = 36.12
